$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New averaging schemes that were run: 3 new "Spiral" rotation schemes.
# They get appended as new rows 17-19 (index values 15, 16, 17 in column A),
# mirroring the structure of the existing rows (A = running index with the
# header style copied from an existing indexed cell, B = scheme name,
# C:M = 1 for each HKL column).

$schemeNames = @(
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space"
)

$startRow = 17
for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $row = $startRow + $i
    $idxValue = 15 + $i

    # Copy the formatting (bold font, border, centered/top alignment) of the
    # existing "A" index column cell down onto the new row, then set its value.
    $ws.Range("A16").Copy($ws.Range("A$row"))
    $ws.Range("A$row").Value = $idxValue

    $ws.Range("B$row").Value = $schemeNames[$i]

    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}
